# Apply cell value updates to existing rows (odds recalculated for 2026-01-07)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 3.55
$ws.Range("Q2").Value = 2.12
$ws.Range("AB2").Value = 6.2
$ws.Range("F4").Value = 1.42
$ws.Range("P4").Value = 3
$ws.Range("AD4").Value = 29
$ws.Range("H5").Value = 3.45
$ws.Range("I5").Value = 3.5
$ws.Range("L5").Value = 1.35
$ws.Range("N5").Value = 4.7
$ws.Range("V5").Value = 1.4
$ws.Range("W5").Value = 1.81
$ws.Range("AF5").Value = 15
$ws.Range("AL5").Value = 32
$ws.Range("P6").Value = 1.8
$ws.Range("Q6").Value = 2.2
$ws.Range("AL6").Value = 44
$ws.Range("H7").Value = 2.42
$ws.Range("R7").Value = 1.37
$ws.Range("AL7").Value = 48
$ws.Range("L8").Value = 1.44
$ws.Range("AN8").Value = 15
$ws.Range("F9").Value = 3.7
$ws.Range("G9").Value = 3.75
$ws.Range("H9").Value = 2.12
$ws.Range("I9").Value = 2.14
$ws.Range("L9").Value = 1.35
$ws.Range("V9").Value = 1.87
$ws.Range("W9").Value = 1.36
$ws.Range("AF9").Value = 27
$ws.Range("AJ9").Value = 70
$ws.Range("S11").Value = 4.7
$ws.Range("S12").Value = 2.74
$ws.Range("AJ12").Value = 330
$ws.Range("U13").Value = 2.12

# Append new row 15: Honduras Liga Nacional - CD Olimpia vs CD Marathon
$ws.Range("A15").Value = "Honduras Liga Nacional"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "2026-01-07"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "22:00:00"
$ws.Range("D15").Value = "CD Olimpia"
$ws.Range("E15").Value = "CD Marathon"
$ws.Range("F15").Value = 1.04
$ws.Range("G15").Value = 1000
$ws.Range("H15").Value = 1.04
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 1.03
$ws.Range("K15").Value = 950
$ws.Range("L15").Value = 1.01
$ws.Range("M15").Value = 1.01
$ws.Range("N15").Value = 1.25
$ws.Range("O15").Value = 1.01
$ws.Range("P15").Value = 1.24
$ws.Range("Q15").Value = 1.01
$ws.Range("R15").Value = 1.18
$ws.Range("S15").Value = 1.41
$ws.Range("T15").Value = 1.04
$ws.Range("U15").Value = 1.04
$ws.Range("V15").Value = 1.01
$ws.Range("W15").Value = 1.01
$ws.Range("X15").Value = 1000
$ws.Range("Y15").Value = 1000
$ws.Range("Z15").Value = 1000
$ws.Range("AA15").Value = 1000
$ws.Range("AB15").Value = 1000
$ws.Range("AC15").Value = 1000
$ws.Range("AD15").Value = 1000
$ws.Range("AE15").Value = 1000
$ws.Range("AF15").Value = 1000
$ws.Range("AG15").Value = 1000
$ws.Range("AH15").Value = 1000
$ws.Range("AI15").Value = 1000
$ws.Range("AJ15").Value = 1000
$ws.Range("AK15").Value = 1000
$ws.Range("AL15").Value = 1000
$ws.Range("AM15").Value = 1000
$ws.Range("AN15").Value = 1000
$ws.Range("AO15").Value = 1000
